$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracking")

# --- Existing "store columns" (H:N) on rows 3-5 gain level-0 flower choices ---
# New shared strings must be introduced in this exact order so they land at
# shared-string indices 65-71 (hibiscus, iris, peony, calla_lily, peace_lily,
# venus_flytrap, daisy_patch) the same way the authoritative edit did.
$ws.Range("I5").Value = "hibiscus"
$ws.Range("H3").Value = "iris"
$ws.Range("L5").Value = "peony"
$ws.Range("J4").Value = "calla_lily"
$ws.Range("K5").Value = "peace_lily"
$ws.Range("M5").Value = "venus_flytrap"
$ws.Range("H4").Value = "daisy_patch"

# --- Row 21 (existing sunplant-less row, now tagged as daisy_patch) ---
$ws.Range("A21").Value = "daisy_patch"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 7
$ws.Range("F21").NumberFormat = "0"
$ws.Range("G21").Value = 1

# --- New tracking rows 22-27 for the level 0 flowers, all drawn by Diconcilio ---
$newRows = @(
    @{ Row = 22; Name = "calla_lily";     D = 16; E = 2 },
    @{ Row = 23; Name = "hibiscus";       D = 16; E = 1 },
    @{ Row = 24; Name = "iris";           D = 16; E = 0 },
    @{ Row = 25; Name = "peace_lily";     D = 16; E = 3 },
    @{ Row = 26; Name = "peony";          D = 16; E = 4 },
    @{ Row = 27; Name = "venus_flytrap";  D = 16; E = 5 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Name
    $ws.Range("B$row").Value = "Diconcilio"
    $ws.Range("C$row").Formula = '=IF(_xlfn.XLOOKUP(B' + $row + ',''Artist Links''!$A$1:$A$5,''Artist Links''!$B$1:$B$5, "")<>"", HYPERLINK(_xlfn.XLOOKUP(B' + $row + ',''Artist Links''!$A$1:$A$5,''Artist Links''!$B$1:$B$5, ""), "Link"), "")'
    $ws.Range("C$row").Style = "Hyperlink"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = 7
    $ws.Range("F$row").NumberFormat = "0"
    $ws.Range("G$row").Value = 1
}

# --- Conditional formatting (color scale) now covers the new rows too ---
$fc = $ws.Range("E2:E20").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("E2:E23"))

# --- Selection left where the author left it ---
$ws.Range("I25").Select()
